$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing rows 2-3 down to 3-4
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# Fill in the new row 2 with the latest news item.
# A2 holds a date-like string ("2026-01-12"); force text format first so
# Excel doesn't auto-convert it to a date serial, then drop the temporary
# number format again so the cell ends up unstyled, like the other cells.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2026-01-12"
$ws.Range("A2").ClearFormats()

$ws.Range("B2").Value = "Madre denuncia a Karina García por grabar a su hija menor de edad sin autorización durante una transmisión en vivo"
$ws.Range("C2").Value = "Infobae"
$ws.Range("D2").Value = "Sin identificar"
$ws.Range("E2").Value = "https://www.infobae.com/colombia/2026/01/12/madre-denuncia-a-karina-garcia-por-grabar-a-su-hija-menor-de-edad-sin-autorizacion-durante-una-transmision-en-vivo/"
$ws.Range("F2").Value = "PorDahana Ospina"
